$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 3.1
$ws.Range("AH2").Value = 7.5
$ws.Range("BB2").Value = 351

$ws.Range("G4").Value = 3.8
$ws.Range("I4").Value = 2.15
$ws.Range("J4").Value = 4.5
$ws.Range("L4").Value = 3
$ws.Range("O4").Value = 1.5
$ws.Range("P4").Value = 2.5
$ws.Range("Q4").Value = 2.6
$ws.Range("R4").Value = 1.48
$ws.Range("S4").Value = 1.57
$ws.Range("T4").Value = 2.25
$ws.Range("W4").Value = 8.5
$ws.Range("AL4").Value = 21
$ws.Range("AO4").Value = 23
$ws.Range("AT4").Value = 2.25
$ws.Range("AZ4").Value = 41

$ws.Range("G5").Value = 1.33
$ws.Range("I5").Value = 12
$ws.Range("L5").Value = 10
$ws.Range("M5").Value = 1.06
$ws.Range("N5").Value = 10
$ws.Range("O5").Value = 1.33
$ws.Range("P5").Value = 3.25
$ws.Range("Q5").Value = 2.08
$ws.Range("R5").Value = 1.73
$ws.Range("U5").Value = 2.63
$ws.Range("V5").Value = 1.44
$ws.Range("AC5").Value = 8
$ws.Range("AD5").Value = 9
$ws.Range("AE5").Value = 29
$ws.Range("AP5").Value = 23
$ws.Range("AS5").Value = 251
$ws.Range("AU5").Value = 12
$ws.Range("AW5").Value = 10
$ws.Range("AY5").Value = 51

$ws.Range("M7").Value = 1.07
$ws.Range("N7").Value = 9
$ws.Range("O7").Value = 1.36
$ws.Range("P7").Value = 3
$ws.Range("Q7").Value = 2.2
$ws.Range("R7").Value = 1.65

$ws.Range("G9").Value = 1.4
$ws.Range("H9").Value = 4.2
$ws.Range("I9").Value = 9
$ws.Range("M9").Value = 1.07
$ws.Range("N9").Value = 9
$ws.Range("Z9").Value = 8.5
$ws.Range("AD9").Value = 8.5
$ws.Range("AW9").Value = 9

$ws.Range("G11").Value = 1.48
$ws.Range("H11").Value = 3.9
$ws.Range("M11").Value = 1.07
$ws.Range("N11").Value = 8.5
$ws.Range("Q11").Value = 2.08
$ws.Range("R11").Value = 1.73
$ws.Range("Z11").Value = 9.5
$ws.Range("AC11").Value = 8.5
$ws.Range("AD11").Value = 8
$ws.Range("AJ11").Value = 21
$ws.Range("AM11").Value = 51

$ws.Range("G14").Value = 7
$ws.Range("H14").Value = 4.5
$ws.Range("I14").Value = 1.44
$ws.Range("M14").Value = 1.05
$ws.Range("N14").Value = 11
$ws.Range("U14").Value = 2.2
$ws.Range("V14").Value = 1.62
$ws.Range("AE14").Value = 23
$ws.Range("AF14").Value = 81
$ws.Range("AI14").Value = 6
$ws.Range("AN14").Value = 8.5
$ws.Range("AR14").Value = 201
$ws.Range("AX14").Value = 7
